$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.183.60"
$ws.Range("E2").Value = "  -2.03%  "
$ws.Range("D3").Value = "1.573.26"
$ws.Range("E3").Value = "  -1.34%  "
$ws.Range("E4").Value = "  -0.38%  "
$ws.Range("D5").Value = "'207.10"
$ws.Range("E5").Value = "  -1.10%  "
$ws.Range("D6").Value = "'0.489"
$ws.Range("E6").Value = "  -2.41%  "
$ws.Range("E7").Value = "  -0.38%  "
$ws.Range("D8").Value = "'22.32"
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").Value = "'0.248"
$ws.Range("E9").Value = "  -1.95%  "
$ws.Range("D10").Value = "'0.0591"
$ws.Range("E10").Value = "  -0.37%  "
$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("D12").Value = "1.795.58"
$ws.Range("E12").Value = "  -1.42%  "
$ws.Range("D13").Value = "1.572.85"
$ws.Range("E13").Value = "  -1.43%  "
$ws.Range("E14").Value = "  -2.05%  "
$ws.Range("D15").Value = "'0.520"
$ws.Range("E15").Value = "  -2.34%  "
$ws.Range("D16").Value = "'62.66"
$ws.Range("E16").Value = "  -1.20%  "
$ws.Range("D17").Value = "27.177.22"
$ws.Range("E17").Value = "  -2.06%  "
$ws.Range("D18").Value = "'215.04"
$ws.Range("E18").Value = "  -1.95%  "
$ws.Range("E19").Value = "  -1.06%  "
$ws.Range("D20").Value = "0.0₃0687"
$ws.Range("E20").Value = "  -1.32%  "
$ws.Range("E21").Value = "  -0.31%  "
$ws.Range("E22").Value = "  -0.41%  "
$ws.Range("D23").Value = "'9.40"
$ws.Range("E23").Value = "  -3.42%  "
$ws.Range("E24").Value = "  +1.11%  "
$ws.Range("D25").Value = "'152.72"
$ws.Range("E25").Value = "  -0.76%  "
$ws.Range("E26").Value = "  -6.59%  "
$ws.Range("D27").Value = "'14.95"
$ws.Range("E27").Value = "  -1.41%  "
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("E30").Value = "  -2.44%  "
$ws.Range("E31").Value = "  -1.94%  "
$ws.Range("E32").Value = "  -1.77%  "
$ws.Range("D33").Value = "1.399.05"
$ws.Range("E33").Value = "  +1.30%  "
$ws.Range("E35").Value = "  +1.13%  "
$ws.Range("E36").Value = "  -2.91%  "
$ws.Range("E37").Value = "  -2.32%  "
$ws.Range("E38").Value = "  -1.66%  "
$ws.Range("D39").Value = "'0.817"
$ws.Range("E39").Value = "  -1.31%  "
$ws.Range("D40").Value = "'0.518"
$ws.Range("E40").Value = "  -3.46%  "
$ws.Range("E41").Value = "  -0.27%  "
$ws.Range("D42").Value = "'0.991"
$ws.Range("E42").Value = "  +1.50%  "
$ws.Range("D43").Value = "'1.81"
$ws.Range("E43").Value = "  +3.33%  "
$ws.Range("D44").Value = "'5.34"
$ws.Range("E44").Value = "  +2.24%  "
$ws.Range("D45").Value = "'63.76"
$ws.Range("E45").Value = "  -1.18%  "
$ws.Range("D46").Value = "'2.18"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").Value = "1.707.57"
$ws.Range("E47").Value = "  -1.42%  "
$ws.Range("D48").Value = "'85.99"
$ws.Range("E48").Value = "  -0.18%  "
$ws.Range("D49").Value = "0.0₇0982"
$ws.Range("E49").Value = "  -2.94%  "
$ws.Range("D50").Value = "'0.0953"
$ws.Range("E50").Value = "  -1.24%  "
$ws.Range("D51").Value = "'0.0494"
$ws.Range("E51").Value = "  -0.36%  "
